# Ibkr.xlsx edit: fix K249 activation flag, append 19 rows of new orders
# (7 MARKET + 11 LIMIT, all sharing the same strike/expiry/target/stop data
# as the row just above them), with the very last row flagged Activation = -1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K249 goes from 1 to 0 (activation already consumed / reset)
$ws.Range("K249").Value = 0

$dateFormat = $ws.Range("E249").NumberFormat

for ($i = 250; $i -le 268; $i++) {
    if ($i -le 256) {
        $orderType = "MARKET"
    } else {
        $orderType = "LIMIT"
    }

    if ($i -eq 268) {
        $activation = -1
    } else {
        $activation = 0
    }

    $ws.Range("A$i").Value = 39250
    $ws.Range("B$i").Value = $orderType
    $ws.Range("C$i").Value = 39400
    $ws.Range("D$i").Value = "PE"
    $ws.Range("E$i").Value = 45660
    $ws.Range("E$i").NumberFormat = $dateFormat
    $ws.Range("F$i").Value = 39500
    $ws.Range("G$i").Value = 39300
    $ws.Range("H$i").Value = 4
    $ws.Range("I$i").Value = 2
    $ws.Range("J$i").Value = 5
    $ws.Range("K$i").Value = $activation
}
